# Add the 2024/10/30 data column (AZ) to the sheet, mirroring the
# existing column layout/styles used by the prior date columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reference cells carrying each of the three cell styles already
# --- used throughout the sheet (plain / yellow highlight / blue highlight).
$styleSrcPlain  = $ws.Range("A2")   # s="1" (メイリオ font, no fill)
$styleSrcYellow = $ws.Range("D2")   # s="2" (yellow fill)
$styleSrcBlue   = $ws.Range("C6")   # s="3" (light-blue fill)
$headerSrc      = $ws.Range("AY1")  # s="1" header/date text style

# --- New column width (matches the other date columns: raw width 12) ---
$ws.Cells.Item(1, 52).EntireColumn.ColumnWidth = 11.17

# --- Header cell AZ1: date text "2024/10/30" ---
# Leading apostrophe forces this to stay literal text instead of being
# parsed/stored as a date serial number (matches the other header cells).
# The value must be assigned *before* the format copy below, otherwise
# the quote-prefix bookkeeping creates/keeps a one-off style instead of
# reusing the shared "s=1" style used by the other header cells.
$ws.Range("AZ1").Value = "'2024/10/30"
$headerSrc.Copy()
$ws.Range("AZ1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows 2-53: (row, styleId, value) ---
$data = @(
    @(2, 2, 115.1),
    @(3, 3, 128.4),
    @(4, 3, 130),
    @(5, 1, 148.8),
    @(6, 1, 157.9),
    @(7, 1, 154),
    @(8, 3, 129.8),
    @(9, 1, 146.7),
    @(10, 1, 175.8),
    @(11, 1, 149.5),
    @(12, 1, 258),
    @(13, 2, 120.2),
    @(14, 1, 151),
    @(15, 1, 148.9),
    @(16, 3, 139.9),
    @(17, 1, 144.1),
    @(18, 1, 157.6),
    @(19, 1, 144.7),
    @(20, 1, 181),
    @(21, 1, 154.5),
    @(22, 1, 180.5),
    @(23, 1, 198.7),
    @(24, 1, 190.9),
    @(25, 1, 189.5),
    @(26, 1, 142),
    @(27, 1, 179.6),
    @(28, 2, 119.1),
    @(29, 1, 155),
    @(30, 1, 171.2),
    @(31, 1, 161.8),
    @(32, 1, 187.6),
    @(33, 1, 156.9),
    @(34, 3, 128.8),
    @(35, 1, 141.2),
    @(36, 2, 123.1),
    @(37, 1, 155.9),
    @(38, 1, 193.3),
    @(39, 1, 172.1),
    @(40, 3, 137.9),
    @(41, 2, 115),
    @(42, 1, 141.8),
    @(43, 1, 165.2),
    @(44, 1, 156.6),
    @(45, 1, 140),
    @(46, 1, 184),
    @(47, 1, 145.9),
    @(48, 1, 241.1),
    @(49, 1, 149.6),
    @(50, 1, 152.8),
    @(51, 1, 188.4),
    @(52, 1, 170.5),
    @(53, 1, 216.8)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $styleId = $entry[1]
    $value = $entry[2]

    if ($styleId -eq 2) {
        $src = $styleSrcYellow
    } elseif ($styleId -eq 3) {
        $src = $styleSrcBlue
    } else {
        $src = $styleSrcPlain
    }

    $src.Copy()
    $ws.Range("AZ$row").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws.Range("AZ$row").Value = $value
}
